$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 55, shifting rows 55-58 down to 56-59
$ws.Rows.Item(55).Insert()

$ws.Range("A55").Value = "Sregeda_wrapper.c"
$ws.Range("B55").Value = "double b_c;"
$ws.Range("C55").Value = "double b_c=0;"

$ws.Range("A53").Select()
